$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginCredentials")

# --- Step 1: move the current row-2 "Customer Interaction" detail block
#     (D2:N2) down to row 4, which currently only has A4:C4 populated. ---
$ws.Range("D2:N2").Copy($ws.Range("D4:N4"))

# --- Step 2: build the new row-3 detail block. Copy the D2:N2 formatting
#     pattern into D3:N3 first so the F/G/H/I "text" number format (and all
#     other styles) match the existing rows, then overwrite with the new
#     values for this ticket/login row. ---
$ws.Range("D2:N2").Copy($ws.Range("D3:N3"))

$ws.Range("D3").Value = 783425592
$ws.Range("E3").Value = "string string string"
$ws.Range("F3").Value = "08-Jul-2020"
$ws.Range("G3").Value = "28 May 2020"
$ws.Range("H3").Value = "01: 17 AM"
$ws.Range("I3").Value = "8925403522001884647"
$ws.Range("J3").Value = "4G"
$ws.Range("K3").Value = 22418813
$ws.Range("L3").Value = 28905441
$ws.Range("M3").Value = "National ID"
$ws.Range("N3").Value = "***22222"

# --- Step 3: row 2 keeps only a new CustomerInteration id in D2; the rest
#     of the old detail block (now duplicated onto row 4) is cleared away
#     completely (formatting included) so the row collapses back to A2:D2. ---
$ws.Range("E2:N2").Clear()
$ws.Range("D2").Value = 782945113

# --- Step 4: selection + page setup housekeeping to match the saved file. ---
$ws.Range("D2").Select()
$ws.PageSetup.Orientation = 1
